$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.82124057048917
$ws.Range("C2").Value = 10.58938781678382
$ws.Range("D2").Value = 14.3139109059934
$ws.Range("E2").Value = 15.29003877717061
$ws.Range("G2").Value = 45.7940489551703
$ws.Range("H2").Value = 18.70076520889824
$ws.Range("J2").Value = 9.153608173510225
$ws.Range("K2").Value = 9.700054605270251
$ws.Range("L2").Value = 11.71034551229484
$ws.Range("M2").Value = 15.91853761412522
$ws.Range("O2").Value = 30.69852010392284

$ws.Range("B3").Value = 13.64745001079453
$ws.Range("C3").Value = 10.59723849102836
$ws.Range("D3").Value = 14.31911197166196
$ws.Range("E3").Value = 15.32134185581483
$ws.Range("G3").Value = 45.92198371876093
$ws.Range("H3").Value = 18.75267659057546
$ws.Range("J3").Value = 9.162557977095073
$ws.Range("K3").Value = 9.563580636293576
$ws.Range("L3").Value = 11.72024039019999
$ws.Range("M3").Value = 15.89829576650625
$ws.Range("O3").Value = 30.78847226914275

$ws.Range("B4").Value = 13.54195430372661
$ws.Range("C4").Value = 10.60253651706377
$ws.Range("D4").Value = 14.32468583258877
$ws.Range("E4").Value = 15.34229776787103
$ws.Range("G4").Value = 46.01014583138709
$ws.Range("H4").Value = 18.78696086766824
$ws.Range("J4").Value = 9.168340564120223
$ws.Range("K4").Value = 9.480163728448975
$ws.Range("L4").Value = 11.7275815944451
$ws.Range("M4").Value = 15.8878324333948
$ws.Range("O4").Value = 30.8486176595603

$ws.Range("B5").Value = 13.49931690667515
$ws.Range("C5").Value = 10.6048160129381
$ws.Range("D5").Value = 14.3275568170128
$ws.Range("E5").Value = 15.35127444635006
$ws.Range("G5").Value = 46.0484845708832
$ws.Range("H5").Value = 18.80153857883077
$ws.Range("J5").Value = 9.170769489008372
$ws.Range("K5").Value = 9.446304404618198
$ws.Range("L5").Value = 11.7308921302233
$ws.Range("M5").Value = 15.8840663124673
$ws.Range("O5").Value = 30.87436274419707

$ws.Range("B6").Value = 13.49225971313097
$ws.Range("C6").Value = 10.60520181132684
$ws.Range("D6").Value = 14.32806978156121
$ws.Range("E6").Value = 15.3527914242183
$ws.Range("G6").Value = 46.05499623000692
$ws.Range("H6").Value = 18.80399584664683
$ws.Range("J6").Value = 9.171177193605756
$ws.Range("K6").Value = 9.440691293743333
$ws.Range("L6").Value = 11.73146112372955
$ws.Range("M6").Value = 15.88347111548307
$ws.Range("O6").Value = 30.87871229536582

$ws.Range("B7").Value = 13.54137778857606
$ws.Range("C7").Value = 10.60256677069227
$ws.Range("D7").Value = 14.32472212275874
$ws.Range("E7").Value = 15.34241706051727
$ws.Range("G7").Value = 46.01065312121706
$ws.Range("H7").Value = 18.78715501111915
$ws.Range("J7").Value = 9.168373027725252
$ws.Range("K7").Value = 9.479706497326786
$ws.Range("L7").Value = 11.72762494927993
$ws.Range("M7").Value = 15.88777962213902
$ws.Range("O7").Value = 30.84895986546882

$ws.Range("B8").Value = 13.76109492851296
$ws.Range("C8").Value = 10.59199582037841
$ws.Range("D8").Value = 14.31521085515881
$ws.Range("E8").Value = 15.30047214998434
$ws.Range("G8").Value = 45.83616377934489
$ws.Range("H8").Value = 18.71816423319913
$ws.Range("J8").Value = 9.156634549660552
$ws.Range("K8").Value = 9.652942874844969
$ws.Range("L8").Value = 11.71349497541246
$ws.Range("M8").Value = 15.9111525639485
$ws.Range("O8").Value = 30.7285153874499

$ws.Range("B9").Value = 14.19945296773748
$ws.Range("C9").Value = 10.5750393937145
$ws.Range("D9").Value = 14.31539470171666
$ws.Range("E9").Value = 15.23196714402258
$ws.Range("G9").Value = 45.57042845864152
$ws.Range("H9").Value = 18.6019791503991
$ws.Range("J9").Value = 9.13588588398107
$ws.Range("K9").Value = 9.993980512505155
$ws.Range("L9").Value = 11.69579919974916
$ws.Range("M9").Value = 15.97241382371303
$ws.Range("O9").Value = 30.53133555388071

$ws.Range("B10").Value = 14.5232391089091
$ws.Range("C10").Value = 10.56485910215129
$ws.Range("D10").Value = 14.32693242111556
$ws.Range("E10").Value = 15.18998660614901
$ws.Range("G10").Value = 45.4220350337384
$ws.Range("H10").Value = 18.528237120475
$ws.Range("J10").Value = 9.122012317632452
$ws.Range("K10").Value = 10.24312445627534
$ws.Range("L10").Value = 11.68886131830639
$ws.Range("M10").Value = 16.02659203052398
$ws.Range("O10").Value = 30.41027328169651

$ws.Range("B11").Value = 14.67031792847373
$ws.Range("C11").Value = 10.56071762733369
$ws.Range("D11").Value = 14.33463583302491
$ws.Range("E11").Value = 15.17269485088537
$ws.Range("G11").Value = 45.3647440159957
$ws.Range("H11").Value = 18.49720707638798
$ws.Range("J11").Value = 9.115995555493342
$ws.Range("K11").Value = 10.35570252303317
$ws.Range("L11").Value = 11.68701215242966
$ws.Range("M11").Value = 16.05317522008836
$ws.Range("O11").Value = 30.36037263817737

$ws.Range("B12").Value = 14.72593370124027
$ws.Range("C12").Value = 10.559219373431
$ws.Range("D12").Value = 14.33790370581032
$ws.Range("E12").Value = 15.16640596411034
$ws.Range("G12").Value = 45.34452128927261
$ws.Range("H12").Value = 18.4858181245688
$ws.Range("J12").Value = 9.113759279606082
$ws.Range("K12").Value = 10.39818706989503
$ws.Range("L12").Value = 11.68649896390733
$ws.Range("M12").Value = 16.06351507204567
$ws.Range("O12").Value = 30.34222055009117

$ws.Range("B13").Value = 14.71396016675182
$ws.Range("C13").Value = 10.55953894036208
$ws.Range("D13").Value = 14.33718434905694
$ws.Range("E13").Value = 15.16774887163131
$ws.Range("G13").Value = 45.34881109178348
$ws.Range("H13").Value = 18.48825486791535
$ws.Range("J13").Value = 9.11423903042142
$ws.Range("K13").Value = 10.38904434902102
$ws.Range("L13").Value = 11.68660118365134
$ws.Range("M13").Value = 16.06127612674561
$ws.Range("O13").Value = 30.34609681946057

$ws.Range("B14").Value = 14.6748953475174
$ws.Range("C14").Value = 10.56059296363436
$ws.Range("D14").Value = 14.33489766684298
$ws.Range("E14").Value = 15.17217226996356
$ws.Range("G14").Value = 45.36305075823816
$ws.Range("H14").Value = 18.49626285803916
$ws.Range("J14").Value = 9.115810732248285
$ws.Range("K14").Value = 10.35920089335012
$ws.Range("L14").Value = 11.68696618854578
$ws.Range("M14").Value = 16.05402043769335
$ws.Range("O14").Value = 30.35886433679838

$ws.Range("B15").Value = 14.65095519062456
$ws.Range("C15").Value = 10.56124769270783
$ws.Range("D15").Value = 14.33354261647438
$ws.Range("E15").Value = 15.174915460263
$ws.Range("G15").Value = 45.37196477486047
$ws.Range("H15").Value = 18.50121504661068
$ws.Range("J15").Value = 9.116778927379121
$ws.Range("K15").Value = 10.34090077701392
$ws.Range("L15").Value = 11.68721409711805
$ws.Range("M15").Value = 16.04961156187272
$ws.Range("O15").Value = 30.36678174381101

$ws.Range("B16").Value = 14.51361864521614
$ws.Range("C16").Value = 10.56513957743284
$ws.Range("D16").Value = 14.32647819598509
$ws.Range("E16").Value = 15.19115295079909
$ws.Range("G16").Value = 45.42598496521948
$ws.Range("H16").Value = 18.53031560226737
$ws.Range("J16").Value = 9.122411434674746
$ws.Range("K16").Value = 10.23574873626825
$ws.Range("L16").Value = 11.68900838971627
$ws.Range("M16").Value = 16.02489329988742
$ws.Range("O16").Value = 30.41363851522064

$ws.Range("B17").Value = 14.42927719558095
$ws.Range("C17").Value = 10.56765225487191
$ws.Range("D17").Value = 14.32277162164217
$ws.Range("E17").Value = 15.20157618789993
$ws.Range("G17").Value = 45.4617431649153
$ws.Range("H17").Value = 18.54881194650743
$ws.Range("J17").Value = 9.125942059381188
$ws.Range("K17").Value = 10.17102015001419
$ws.Range("L17").Value = 11.69044324865114
$ws.Range("M17").Value = 16.01022199836193
$ws.Range("O17").Value = 30.44370852703789

$ws.Range("B18").Value = 14.38074903439907
$ws.Range("C18").Value = 10.56914358508142
$ws.Range("D18").Value = 14.32087090857372
$ws.Range("E18").Value = 15.20774132115295
$ws.Range("G18").Value = 45.48327153743527
$ws.Range("H18").Value = 18.55968734909595
$ws.Range("J18").Value = 9.128000502340639
$ws.Range("K18").Value = 10.13372071116343
$ws.Range("L18").Value = 11.69139160988266
$ws.Range("M18").Value = 16.001966123176
$ws.Range("O18").Value = 30.46149070664611

$ws.Range("B19").Value = 14.36431679750977
$ws.Range("C19").Value = 10.56965645367548
$ws.Range("D19").Value = 14.32026713512172
$ws.Range("E19").Value = 15.20985793534607
$ws.Range("G19").Value = 45.49072565200102
$ws.Range("H19").Value = 18.56341025197107
$ws.Range("J19").Value = 9.128702222745707
$ws.Range("K19").Value = 10.12108101988972
$ws.Range("L19").Value = 11.69173387191491
$ws.Range("M19").Value = 15.99920235393563
$ws.Range("O19").Value = 30.46759502033989

$ws.Range("B20").Value = 14.43825762584816
$ws.Range("C20").Value = 10.56738000710188
$ws.Range("D20").Value = 14.3231422789499
$ws.Range("E20").Value = 15.2004490294958
$ws.Range("G20").Value = 45.45783713687857
$ws.Range("H20").Value = 18.54681847379682
$ws.Range("J20").Value = 9.125563350729509
$ws.Range("K20").Value = 10.17791806114403
$ws.Range("L20").Value = 11.69027777460154
$ws.Range("M20").Value = 16.0117649147821
$ws.Range("O20").Value = 30.44045714646385

$ws.Range("B21").Value = 14.6863721825143
$ws.Range("C21").Value = 10.56028147371927
$ws.Range("D21").Value = 14.33555982039227
$ws.Range("E21").Value = 15.1708659820083
$ws.Range("G21").Value = 45.35882824455273
$ws.Range("H21").Value = 18.49390090975735
$ws.Range("J21").Value = 9.115347943140122
$ws.Range("K21").Value = 10.36797090100604
$ws.Range("L21").Value = 11.68685390852089
$ws.Range("M21").Value = 16.05614423271345
$ws.Range("O21").Value = 30.35509400551526

$ws.Range("B22").Value = 14.84804249005375
$ws.Range("C22").Value = 10.55605020584553
$ws.Range("D22").Value = 14.34571868555274
$ws.Range("E22").Value = 15.15304189765047
$ws.Range("G22").Value = 45.30270224797719
$ws.Range("H22").Value = 18.4614229142797
$ws.Range("J22").Value = 9.10891712794316
$ws.Range("K22").Value = 10.49131221595245
$ws.Range("L22").Value = 11.68570604561739
$ws.Range("M22").Value = 16.0867395477356
$ws.Range("O22").Value = 30.30364228484456

$ws.Range("B23").Value = 14.76181645194971
$ws.Range("C23").Value = 10.55827129971111
$ws.Range("D23").Value = 14.34011054778744
$ws.Range("E23").Value = 15.16241693297063
$ws.Range("G23").Value = 45.33187147449256
$ws.Range("H23").Value = 18.47856436239119
$ws.Range("J23").Value = 9.1123269707139
$ws.Range("K23").Value = 10.42557411883345
$ws.Range("L23").Value = 11.68621926525037
$ws.Range("M23").Value = 16.07026646511828
$ws.Range("O23").Value = 30.33070592234291

$ws.Range("B24").Value = 14.43419769147421
$ws.Range("C24").Value = 10.56750294464415
$ws.Range("D24").Value = 14.32297398743182
$ws.Range("E24").Value = 15.20095807991066
$ws.Range("G24").Value = 45.45960002933214
$ws.Range("H24").Value = 18.54771897018441
$ws.Range("J24").Value = 9.125734475685679
$ws.Range("K24").Value = 10.17479978089695
$ws.Range("L24").Value = 11.69035220084788
$ws.Range("M24").Value = 16.01106680477484
$ws.Range("O24").Value = 30.44192555522505

$ws.Range("B25").Value = 14.08035824877605
$ws.Range("C25").Value = 10.57922495894337
$ws.Range("D25").Value = 14.31333588948667
$ws.Range("E25").Value = 15.24903087957948
$ws.Range("G25").Value = 45.63410721524078
$ws.Range("H25").Value = 18.63136783253325
$ws.Range("J25").Value = 9.141257297527291
$ws.Range("K25").Value = 9.901811325287795
$ws.Range("L25").Value = 11.69951847521005
$ws.Range("M25").Value = 15.9542127547877
$ws.Range("O25").Value = 30.58049894665852
